$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 3")

# Update selection to F9 (matches diff <selection activeCell="F9" sqref="F9"/>)
$ws.Range("F9").Select()

# Row 9: set Interuption Time (E9) and update Delta Time (F9)
$ws.Range("E9").Value = 20
$ws.Range("F9").Value = 30

# Row 10: fill in Date, Start, Activity and Comments
$ws.Range("B10").Value = 43878
$ws.Range("C10").Value = 0.77083333333333337
$ws.Range("G10").Value = "kood"
$ws.Range("H10").Value = "HW2 lõpuni tegemine"

$wb.Save()
